$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as text so numeric-looking values
# (e.g. "298.87") are stored as strings, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.155.09'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.265.27'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '298.87'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").Value = '95.37'
$ws.Range("E6").Value = '  -5.16%  '
$ws.Range("E7").Value = '  -2.47%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.74%  '
$ws.Range("D10").Value = '33.27'
$ws.Range("E10").Value = '  -3.72%  '
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").Value = '47.85'
$ws.Range("E12").Value = '  -8.27%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '6.66'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = '2.619.24'
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").Value = '15.48'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").Value = '2.282.34'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '0.782'
$ws.Range("E18").Value = '  -5.27%  '
$ws.Range("D19").Value = '42.078.85'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").Value = '11.67'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").Value = '66.47'
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").Value = '234.08'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("D28").Value = '23.92'
$ws.Range("E28").Value = '  -6.65%  '
$ws.Range("D29").Value = '2.27'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").Value = '168.14'
$ws.Range("E30").Value = '  +4.92%  '
$ws.Range("D31").Value = '9.15'
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").Value = '33.52'
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("D35").Value = '4.46'
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("E36").Value = '  -4.86%  '
$ws.Range("D37").Value = '16.49'
$ws.Range("E37").Value = '  -3.13%  '
$ws.Range("D38").Value = '0.0684'
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("E39").Value = '  -4.08%  '
$ws.Range("D40").Value = '0.0982'
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("E41").Value = '  -2.85%  '
$ws.Range("E42").Value = '  -5.84%  '
$ws.Range("D43").Value = '2.44'
$ws.Range("E43").Value = '  -4.53%  '
$ws.Range("D44").Value = '1.957.70'
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("D45").Value = '0.0276'
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D46").Value = '17.35'
$ws.Range("E46").Value = '  -7.64%  '
$ws.Range("D47").Value = '9.52'
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("E48").Value = '  -4.01%  '
$ws.Range("D49").Value = '2.491.04'
$ws.Range("E49").Value = '  -2.33%  '
$ws.Range("D50").Value = '52.14'
$ws.Range("E50").Value = '  -6.31%  '
$ws.Range("D51").Value = '4.51'
$ws.Range("E51").Value = '  -3.70%  '

# Restore original (default) cell style on the Price column.
$ws.Range("D2:D51").Style = "Normal"
